$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Item" bullet paragraph: fix wording and emphasize "can"
#    " ... Items get an inherent -1B, but can't be used ..."
#    -> " ... Items get an -1 Burn, but can't be used ..."
#    and the word "can" right after "... and Sword " becomes bold+italic
# ---------------------------------------------------------------------
$body = $d.Content
$res1 = $body.Find.Execute("inherent -1B", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "-1 Burn", 2)
Write-Output "fix inherent-1B: $res1"

$body2 = $d.Content
$res2 = $body2.Find.Execute("Sword can be designated", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
Write-Output "find Sword-can: $res2 text=[$($body2.Text)]"
if ($res2) {
    $canStart = $body2.Start + 6
    $canEnd = $canStart + 3
    $canRng = $d.Range($canStart, $canEnd)
    Write-Output "canRng=[$($canRng.Text)]"
    $canRng.Bold = 1
    $canRng.Italic = 1
}

# ---------------------------------------------------------------------
# 2) Powers table (2nd table) numeric / text tweaks
# ---------------------------------------------------------------------
$t = $d.Tables.Item(2)

# Bow row (2): Burn 4 -> 5 ; Effects Pierce(2)->Pierce (2), Knock(1)->Knock (1)
$t.Cell(2, 8).Range.Text = "5"
$bowFx = $t.Cell(2, 9).Range
$null = $bowFx.Find.Execute("Pierce(", $true, $false, $false, $false, $false, $true, 1, $false, "Pierce (", 2)
$bowFx2 = $t.Cell(2, 9).Range
$null = $bowFx2.Find.Execute("Knock(", $true, $false, $false, $false, $false, $true, 1, $false, "Knock (", 2)

# Breast Plate row (3): Enhancements "Resist Piece(1)" -> "Resist Pierce (1)"
$bpEnh = $t.Cell(3, 10).Range
$null = $bpEnh.Find.Execute("Piece(", $true, $false, $false, $false, $false, $true, 1, $false, "Pierce (", 2)

# Horn row (9): Burn 2 -> 3
$t.Cell(9, 8).Range.Text = "3"

# Net row (11): Burn 4 -> 5 ; Effects "4/4/-" -> "4/2/-"
$t.Cell(11, 8).Range.Text = "5"
$netFx = $t.Cell(11, 9).Range
$null = $netFx.Find.Execute("4/4/", $true, $false, $false, $false, $false, $true, 1, $false, "4/2/", 2)

# Shield row (13): Burn 2 -> 3 ; Effects "Block(" -> "Block ("
$t.Cell(13, 8).Range.Text = "3"
$shFx = $t.Cell(13, 9).Range
$null = $shFx.Find.Execute("Block(", $true, $false, $false, $false, $false, $true, 1, $false, "Block (", 2)

# Sword row (15): Burn 3 -> 4
$t.Cell(15, 8).Range.Text = "4"

Write-Output "done"
